$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.423.26'
$ws.Range("E2").Value = '  -1.68%  '

# Row 3
$ws.Range("D3").Value = '3.541.53'
$ws.Range("E3").Value = '  -2.96%  '

# Row 4
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
$ws.Range("D5").Value = '584.11'
$ws.Range("E5").Value = '  +1.17%  '

# Row 6
$ws.Range("D6").Value = '172.96'
$ws.Range("E6").Value = '  -1.85%  '

# Row 7
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").Value = '0.613'
$ws.Range("E7").Value = '  +0.06%  '

# Row 8
$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D8").Value = '3.535.43'
$ws.Range("E8").Value = '  -2.92%  '

# Row 9
$ws.Range("E9").Value = '  +0.09%  '

# Row 10
$ws.Range("D10").Value = '0.191'
$ws.Range("E10").Value = '  -3.36%  '

# Row 11
$ws.Range("D11").Value = '6.80'
$ws.Range("E11").Value = '  -0.65%  '

# Row 12
$ws.Range("D12").Value = '0.584'
$ws.Range("E12").Value = '  -3.41%  '

# Row 13
$ws.Range("D13").Value = '47.52'
$ws.Range("E13").Value = '  -2.27%  '

# Row 14
$ws.Range("E14").Value = '  -3.85%  '

# Row 15
$ws.Range("D15").Value = '4.112.17'

# Row 16
$ws.Range("D16").Value = '8.55'
$ws.Range("E16").Value = '  -3.82%  '

# Row 17
$ws.Range("D17").Value = '630.45'
$ws.Range("E17").Value = '  -6.07%  '

# Row 18
$ws.Range("D18").Value = '3.539.56'
$ws.Range("E18").Value = '  -3.08%  '

# Row 19
$ws.Range("D19").Value = '69.468.54'
$ws.Range("E19").Value = '  -1.71%  '

# Row 20
$ws.Range("E20").Value = '  +1.38%  '

# Row 21
$ws.Range("D21").Value = '17.43'
$ws.Range("E21").Value = '  -1.99%  '

# Row 22
$ws.Range("D22").Value = '11.21'

# Row 23
$ws.Range("D23").Value = '0.892'

# Row 24
$ws.Range("D24").Value = '15.99'
$ws.Range("E24").Value = '  -6.54%  '

# Row 25
$ws.Range("D25").Value = '97.45'
$ws.Range("E25").Value = '  -3.08%  '

# Row 26
$ws.Range("D26").Value = '3.81'
$ws.Range("E26").Value = '  -2.19%  '

# Row 27
$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.12%  '

# Row 28
$ws.Range("B28").Value = 'ImmutableX'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D28").Value = '2.65'
$ws.Range("E28").Value = '  -5.01%  '

# Row 29
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = '9.36'
$ws.Range("E29").Value = '  -6.25%  '

# Row 30
$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").Value = '32.94'
$ws.Range("E30").Value = '  -5.79%  '

# Row 31
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = '8.60'
$ws.Range("E31").Value = '  -4.35%  '

# Row 32
$ws.Range("B32").Value = 'Stacks'
$ws.Range("C32").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D32").Value = '3.16'
$ws.Range("E32").Value = '  -5.79%  '

# Row 33
$ws.Range("B33").Value = 'Mantle'
$ws.Range("C33").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D33").Value = '1.34'
$ws.Range("E33").Value = '  -3.58%  '

# Row 34
$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").Value = '7.02'
$ws.Range("E34").Value = '  -3.97%  '

# Row 35
$ws.Range("B35").Value = 'Bittensor'
$ws.Range("C35").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D35").Value = '638.49'
$ws.Range("E35").Value = '  +9.12%  '

# Row 36
$ws.Range("B36").Value = 'Cosmos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D36").Value = '10.82'
$ws.Range("E36").Value = '  -2.08%  '

# Row 37
$ws.Range("B37").Value = 'dogwifhat'
$ws.Range("C37").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D37").Value = '3.53'
$ws.Range("E37").Value = '  -11.52%  '

# Row 38
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = '0.103'
$ws.Range("E38").Value = '  -3.44%  '

# Row 39
$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").Value = '57.32'
$ws.Range("E39").Value = '  -1.39%  '

# Row 40
$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  +0.03%  '

# Row 41
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.0457'
$ws.Range("E41").Value = '  +0.82%  '

# Row 42
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").Value = '0.137'
$ws.Range("E42").Value = '  -3.50%  '

# Row 43
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '3.393.55'
$ws.Range("E43").Value = '  -5.21%  '

# Row 44
$ws.Range("B44").Value = 'TheGraph'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D44").Value = '0.330'
$ws.Range("E44").Value = '  -4.14%  '

# Row 45
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = '32.90'
$ws.Range("E45").Value = '  -5.44%  '

# Row 46
$ws.Range("E46").Value = '  -5.35%  '

# Row 47
$ws.Range("B47").Value = 'Fetch.AI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D47").Value = '2.57'
$ws.Range("E47").Value = '  -4.89%  '

# Row 48
$ws.Range("B48").Value = 'ThetaToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D48").Value = '2.75'
$ws.Range("E48").Value = '  -5.94%  '

# Row 49
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").Value = '0.130'
$ws.Range("E49").Value = '  -2.17%  '

# Row 50
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D50").Value = '132.46'
$ws.Range("E50").Value = '  -2.35%  '

# Row 51
$ws.Range("B51").Value = 'MXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D51").Value = '5.67'
$ws.Range("E51").Value = '  +13.70%  '
